$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "estimate" label in A2 from "denominator" to "co2_offset_per_MWh"
# now that the subsidy payment's denominator has been properly isolated.
$ws.Range("A2").Value = "co2_offset_per_MWh"

# Move the active selection to A3 (next empty row) as left by the author.
$ws.Range("A3").Select()
